$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New link strings (used for both the displayed cell text and the hyperlink target)
$linkGenomeweb = "https://www.genomeweb.com/regulatory-news-fda-approvals/fda-expands-labels-roches-her2-cdx-assays-id-breast-cancer-patients"
$link360dx     = "https://www.360dx.com/regulatory-news-fda-approvals/fda-expands-labels-roches-her2-cdx-assays-id-breast-cancer-patients"

# New title string shared by both new rows
$title = "FDA Expands Labels For Roche's HER2 CDx Assays to ID Breast Cancer Patients Eligible For Enhertu"

# Keywords value used for both new rows (same as existing rows using "CDx")
$keywords = "CDx"

# Reference style (Hyperlink cell style) taken from the last existing data row's link cell
$linkStyle = $ws.Range("A71").Style

# Row 72
$ws.Range("A72").Value = $linkGenomeweb
$ws.Range("B72").Value = $keywords
$ws.Range("C72").Value = $title

# Row 73
$ws.Range("A73").Value = $link360dx
$ws.Range("B73").Value = $keywords
$ws.Range("C73").Value = $title

# Add the hyperlinks (this also stamps a hyperlink-like style on the cells),
# then reapply the workbook's existing "Hyperlink" cell style so the new
# link cells are formatted exactly like the rest of column A.
$ws.Hyperlinks.Add($ws.Range("A72"), $linkGenomeweb) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A73"), $link360dx) | Out-Null
$ws.Range("A72:A73").Style = $linkStyle

Write-Host "Added rows 72 and 73"
